$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit cyclically rotates the data rows 4-10 (each row now holds the
# values that used to live 3 rows below it, wrapping around within 4-10).
# Row 4 <- old row 7
# Row 5 <- old row 8
# Row 6 <- old row 9
# Row 7 <- old row 10
# Row 8 <- old row 4
# Row 9 <- old row 5
# Row 10 <- old row 6

$ws.Range("A4").Value = 111910419
$ws.Range("B4").Value = 89423
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 5432
$ws.Range("F4").Value = "Granticka"
$ws.Range("G4").Value = "Porodaedalea chrysoloma"
$ws.Range("H4").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("M4").ClearContents()
$ws.Range("Q4").Value = 408379.6130448866
$ws.Range("R4").Value = 7020248.537071504
$ws.Range("AC4").ClearContents()

$ws.Range("A5").Value = 111910415
$ws.Range("B5").Value = 78612
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 6464
$ws.Range("F5").Value = "Luddlav"
$ws.Range("G5").Value = "Nephroma resupinatum"
$ws.Range("H5").Value = "(L.) Ach."
$ws.Range("Q5").Value = 408381.4174405072
$ws.Range("R5").Value = 7020090.541921036

$ws.Range("A6").Value = 111910416
$ws.Range("B6").Value = 73510
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 6428
$ws.Range("F6").Value = "Rostfläck"
$ws.Range("G6").Value = "Arthonia vinosa"
$ws.Range("H6").Value = "Leight."
$ws.Range("Q6").Value = 408382.2393677595
$ws.Range("R6").Value = 7020103.568264721
$ws.Range("AJ6").Value = "gran"
$ws.Range("AK6").Value = "Picea abies"
$ws.Range("AO6").Value = "Picea abies"

$ws.Range("A7").Value = 111910403
$ws.Range("B7").Value = 77597
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 864
$ws.Range("F7").Value = "Knottrig blåslav"
$ws.Range("G7").Value = "Hypogymnia bitteri"
$ws.Range("H7").Value = "(Lynge) Ahti"
$ws.Range("Q7").Value = 408251.2666143124
$ws.Range("R7").Value = 7019793.194737672

$ws.Range("A8").Value = 111910567
$ws.Range("B8").Value = 56398
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 100109
$ws.Range("F8").Value = "Tretåig hackspett"
$ws.Range("G8").Value = "Picoides tridactylus"
$ws.Range("H8").Value = "(Linnaeus, 1758)"
$ws.Range("M8").Value = "färska spår"
$ws.Range("Q8").Value = 408284.1964350128
$ws.Range("R8").Value = 7019857.509490959
$ws.Range("AC8").Value = "ringhack"

$ws.Range("A9").Value = 111910584
$ws.Range("B9").Value = 89423
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 5432
$ws.Range("F9").Value = "Granticka"
$ws.Range("G9").Value = "Porodaedalea chrysoloma"
$ws.Range("H9").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q9").Value = 408422.6472976486
$ws.Range("R9").Value = 7020304.006376172
$ws.Range("AJ9").ClearContents()
$ws.Range("AK9").ClearContents()
$ws.Range("AO9").ClearContents()

$ws.Range("A10").Value = 111910425
$ws.Range("B10").Value = 76499
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 228579
$ws.Range("F10").Value = "Liten svartspik"
$ws.Range("G10").Value = "Chaenothecopsis nana"
$ws.Range("H10").Value = "Tibell"
$ws.Range("Q10").Value = 408400.3469668561
$ws.Range("R10").Value = 7020264.594941486
